# Add payment 71277620 (Cash) 2025-08-18T17:04:07
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing row 28's phone number cell was stored as text in the source
# file; the new append normalizes it to a real number, matching the
# numeric storage used by the rest of that phone number's rows (21-27).
$ws.Range("A28").Value = 71277628

# New row 29 - a Cash payment for phone 71277620.
# The phone number is stored as text (matches the source export's typing),
# so force a quote-prefixed text entry and then strip the quote-prefix
# style back to Normal so no extra cell style sticks around.
$ws.Range("A29").Value = "'71277620"
$ws.Range("A29").Style = "Normal"

# amount column is blank for this payment (blank/text like the rows above it)
$ws.Range("B29").Value = "'"
$ws.Range("B29").Style = "Normal"

$ws.Range("C29").Value = "Cash"
$ws.Range("D29").Value = "2025-08-18T17:04:07"
$ws.Range("E29").Value = 76

# discount_applied column is blank for this payment
$ws.Range("F29").Value = "'"
$ws.Range("F29").Style = "Normal"

$ws.Range("G29").Value = 76
$ws.Range("H29").Value = 0
$ws.Range("I29").Value = 0
